$d = $word.ActiveDocument

# --- Replace "<id>113r_1</id>" with "<id>p113r_1</id>" ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "<id>113r_1</id>", $false, $false, $false, $false, $false, $true, 1, $false,
    "<id>p113r_1</id>", 2)

# --- Replace "<id>113v_1</id>" with "<id>p113v_1</id>" ---
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(
    "<id>113v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false,
    "<id>p113v_1</id>", 2)
